$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33, shifting the existing data (old rows 33-121)
# down to rows 34-122. Excel carries the row formatting (e.g. the date style
# on column D) down with the shifted cells automatically.
$ws.Rows("33").Insert()

# Populate the newly inserted row 33 with the new weekly record.
$ws.Range("A33").Value = 7
$ws.Range("B33").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C33").Value = "Ñuble"
$ws.Range("D33").Value = 44459
$ws.Range("E33").Value = 16
$ws.Range("F33").Value = 100112003
$ws.Range("G33").Value = "Ajo"
$ws.Range("H33").Value = "Chino"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 60
$ws.Range("K33").Value = 16000
$ws.Range("L33").Value = 17000
$ws.Range("M33").Value = 16500
$ws.Range("N33").Value = "$/caja 10 kilos"
$ws.Range("O33").Value = "China"
$ws.Range("P33").Value = 1650
$ws.Range("Q33").Value = 10
$ws.Range("R33").Value = "Hortaliza"
